$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write column D (session directory) first for all new rows, then A, then B, to match
# the shared-string creation order expected by the target workbook.
for ($r = 77; $r -le 91; $r++) {
  $ws.Cells.Item($r, 4).Value = "2013-03-25-anderson"
}
for ($r = 77; $r -le 91; $r++) {
  $ws.Cells.Item($r, 1).Value = "S6"
}
for ($r = 77; $r -le 91; $r++) {
  $ws.Cells.Item($r, 2).Value = "Anderson"
}

$filenames = @(
  "2013-03-25-15-22-16",
  "2013-03-25-15-38-39",
  "2013-03-25-16-01-53",
  "2013-03-25-16-09-13",
  "2013-03-25-16-16-03",
  "2013-03-25-16-25-57",
  "2013-03-25-16-32-17",
  "2013-03-25-16-38-06",
  "2013-03-25-16-44-42",
  "2013-03-25-16-50-37",
  "2013-03-25-17-01-48",
  "2013-03-25-17-07-53",
  "2013-03-25-17-13-59",
  "2013-03-25-17-21-29",
  "2013-03-25-17-27-31"
)
$conditions = @(
  "hybrid-8-57Hz","hybrid-10Hz","hybrid-12Hz","hybrid-15Hz","oddball",
  "oddball","hybrid-10Hz","hybrid-12Hz","hybrid-15Hz","hybrid-8-57Hz",
  "hybrid-12Hz","hybrid-15Hz","hybrid-8-57Hz","oddball","hybrid-10Hz"
)
$runs = @(1,1,1,1,1,2,2,2,2,2,3,3,3,3,3)

for ($i = 0; $i -lt 15; $i++) {
  $r = 77 + $i
  $ws.Cells.Item($r, 5).Value = $filenames[$i]
  $ws.Cells.Item($r, 6).Value = $conditions[$i]
  $ws.Cells.Item($r, 7).Value = $runs[$i]
  $ws.Cells.Item($r, 3).Value = 41358
}

# Match style/format of column C (date) to the existing date cells
$ws.Range("C76").Copy()
$ws.Range("C77:C91").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column D grew a bit wider to fit the new, longer session-directory name
$ws.Columns.Item(4).ColumnWidth = 18.86

# Scroll the view down to the newly-added rows and select the new corner cell
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("J91").Select()
